# Update cryptos list (prices/volume%) per commit
# Numeric-looking Price (column D) values are prefixed with a leading
# apostrophe so Excel stores them as text (matching the source data's
# inline-string cells) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.894.20'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '1.641.71'
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = '  +0.74%  '
$ws.Range("D5").Value = "'215.25"
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").Value = "'0.5066"
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("D8").Value = "'0.2574"
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").Value = "'0.06408"
$ws.Range("D10").Value = "'19.73"
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("D11").Value = "'0.07779"
$ws.Range("E11").Value = '  +1.62%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.766.02'
$ws.Range("E12").Value = '  +6.57%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'4.306"
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("D14").Value = "'0.5451"
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("D15").Value = '0.0₅7899'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").Value = "'65.18"
$ws.Range("E16").Value = '  +2.72%  '
$ws.Range("D17").Value = '25.993.32'
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").Value = "'1.007"
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("D19").Value = "'197.40"
$ws.Range("E19").Value = '  -2.45%  '
$ws.Range("D20").Value = "'4.409"
$ws.Range("E20").Value = '  +2.44%  '
$ws.Range("D21").Value = "'9.989"
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("D22").Value = "'6.050"
$ws.Range("E22").Value = '  +1.75%  '
$ws.Range("D23").Value = "'1.011"
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("D24").Value = "'1.860"
$ws.Range("E24").Value = '  -3.76%  '
$ws.Range("D25").Value = "'140.51"
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").Value = "'0.1148"
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").Value = "'6.895"
$ws.Range("E27").Value = '  +3.21%  '
$ws.Range("D28").Value = "'15.71"
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = "'0.05065"
$ws.Range("E29").Value = '  +1.92%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'1.243"
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").Value = "'3.269"
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").Value = "'3.197"
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("D33").Value = "'1.538"
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("D34").Value = "'2.369"
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("D35").Value = "'0.8940"
$ws.Range("E35").Value = '  +0.29%  '
$ws.Range("D36").Value = "'2.600"
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("D37").Value = "'0.5538"
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("D38").Value = '1.130.87'
$ws.Range("E38").Value = '  -3.51%  '
$ws.Range("D39").Value = "'0.01566"
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("D41").Value = "'5.677"
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("D42").Value = "'0.8160"
$ws.Range("E42").Value = '  +1.74%  '
$ws.Range("D43").Value = "'99.74"
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("E44").Value = '  +7.79%  '
$ws.Range("D45").Value = '1.780.25'
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").Value = "'0.4540"
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = "'1.007"
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = "'55.19"
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").Value = "'0.05092"
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("E50").Value = '  +0.60%  '
$ws.Range("D51").Value = "'0.09566"
$ws.Range("E51").Value = '  +3.19%  '
